# Apply updated cryptocurrency price / volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '28.569.90'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  +0.66%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.565.85'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -1.32%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.10%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '211.94'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  -1.07%  '; ForceText = $false },
    @{ Cell = 'E6'; Value = '  -0.59%  '; ForceText = $false },
    @{ Cell = 'E7'; Value = '  -0.08%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '46.37'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  +3.38%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '24.20'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  +1.71%  '; ForceText = $false },
    @{ Cell = 'E10'; Value = '  -1.58%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.0592'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -1.33%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  -0.44%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '1.786.97'; ForceText = $false },
    @{ Cell = 'E13'; Value = '  -1.44%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '1.561.44'; ForceText = $false },
    @{ Cell = 'E14'; Value = '  -1.65%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '0.522'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  -1.80%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '28.550.49'; ForceText = $false },
    @{ Cell = 'E16'; Value = '  +0.55%  '; ForceText = $false },
    @{ Cell = 'D17'; Value = '3.68'; ForceText = $true },
    @{ Cell = 'E17'; Value = '  -2.90%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '62.13'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -1.80%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '228.61'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  -1.32%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '0.0₃0695'; ForceText = $false },
    @{ Cell = 'E20'; Value = '  -1.90%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '7.34'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -1.87%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  -0.12%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '3.87'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -5.97%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '9.14'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  -2.66%  '; ForceText = $false },
    @{ Cell = 'E25'; Value = '  +6.50%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '151.00'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -0.68%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '14.98'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -1.76%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '6.45'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -2.49%  '; ForceText = $false },
    @{ Cell = 'E29'; Value = '  -3.44%  '; ForceText = $false },
    @{ Cell = 'E30'; Value = '  -0.07%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '0.0466'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -1.28%  '; ForceText = $false },
    @{ Cell = 'E32'; Value = '  -3.51%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '3.21'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -1.01%  '; ForceText = $false },
    @{ Cell = 'E34'; Value = '  -0.74%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '1.398.90'; ForceText = $false },
    @{ Cell = 'E35'; Value = '  -0.57%  '; ForceText = $false },
    @{ Cell = 'E36'; Value = '  -0.92%  '; ForceText = $false },
    @{ Cell = 'E37'; Value = '  -2.92%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '2.35'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +0.75%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '2.57'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  +1.52%  '; ForceText = $false },
    @{ Cell = 'E40'; Value = '  -1.06%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.536'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -1.18%  '; ForceText = $false },
    @{ Cell = 'E42'; Value = '  -0.07%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.789'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  -3.37%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '1.88'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  +2.88%  '; ForceText = $false },
    @{ Cell = 'E45'; Value = '  -4.38%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '0.977'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -0.54%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '62.74'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -2.79%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '1.700.64'; ForceText = $false },
    @{ Cell = 'E48'; Value = '  -1.57%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '86.42'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -1.39%  '; ForceText = $false },
    @{ Cell = 'E50'; Value = '  -4.76%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '0.0518'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -0.79%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Source values are plain-looking numbers (e.g. "46.37") that must stay
        # text, matching the original inline-string cells (avoids "151.00" -> 151).
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"
